$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.384145666666667
$ws.Range("H2").Value = 4.152437
$ws.Range("I2").Value = 0.1014617184198512
$ws.Range("J2").Value = 0.1334061399754118
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.495057333333333
$ws.Range("N2").Value = 16.485172
$ws.Range("O2").Value = 0.8161989011161211
$ws.Range("P2").Value = 0.8403205285996808
$ws.Range("Q2").Value = 7.605959796018221
$ws.Range("R2").Value = 68.45363816416399
$ws.Range("S2").Value = 0.08281294307963584
$ws.Range("T2").Value = 0.112103918062581

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.384145666666667
$ws.Range("H3").Value = 4.152437
$ws.Range("I3").Value = 0.1014617184198512
$ws.Range("J3").Value = 0.1334061399754118
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.09768528951377062
$ws.Range("P3").Value = 0.1005722428790014
$ws.Range("Q3").Value = 0.910305544014
$ws.Range("R3").Value = 8.192749896126
$ws.Range("S3").Value = 0.009911317338407836
$ws.Range("T3").Value = 0.01341695471115717

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.384145666666667
$ws.Range("H4").Value = 4.152437
$ws.Range("I4").Value = 0.1014617184198512
$ws.Range("J4").Value = 0.1334061399754118
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5797745
$ws.Range("N4").Value = 1.159549
$ws.Range("O4").Value = 0.08611580937010824
$ws.Range("P4").Value = 0.0591072285213179
$ws.Range("Q4").Value = 0.8024923618188333
$ws.Range("R4").Value = 4.814954170912999
$ws.Range("S4").Value = 0.008737458001807504
$ws.Range("T4").Value = 0.007885267201673587

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.458038666666667
$ws.Range("H5").Value = 7.374116000000001
$ws.Range("I5").Value = 0.1801810554109116
$ws.Range("J5").Value = 0.2369096391566985
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.495057333333333
$ws.Range("N5").Value = 16.485172
$ws.Range("O5").Value = 0.8161989011161211
$ws.Range("P5").Value = 0.8403205285996808
$ws.Range("Q5").Value = 13.50706340088355
$ws.Range("R5").Value = 121.563570607952
$ws.Range("S5").Value = 0.147063579428329
$ws.Range("T5").Value = 0.1990800332065165

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.458038666666667
$ws.Range("H6").Value = 7.374116000000001
$ws.Range("I6").Value = 0.1801810554109116
$ws.Range("J6").Value = 0.2369096391566985
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.657666
$ws.Range("N6").Value = 1.972998
$ws.Range("O6").Value = 0.09768528951377062
$ws.Range("P6").Value = 0.1005722428790014
$ws.Range("Q6").Value = 1.616568457752
$ws.Range("R6").Value = 14.549116119768
$ws.Range("S6").Value = 0.01760103856271164
$ws.Range("T6").Value = 0.02382653376964406

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.1801810554109116
$ws.Range("J7").Value = 0.2369096391566985
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5797745
$ws.Range("N7").Value = 1.159549
$ws.Range("O7").Value = 0.08611580937010824
$ws.Range("P7").Value = 0.0591072285213179
$ws.Range("Q7").Value = 1.425108138947333
$ws.Range("R7").Value = 8.550648833684
$ws.Range("S7").Value = 0.01551643741987097
$ws.Range("T7").Value = 0.01400307218053794

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 9.799864
$ws.Range("H8").Value = 19.599728
$ws.Range("I8").Value = 0.7183572261692373
$ws.Range("J8").Value = 0.6296842208678898
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.495057333333333
$ws.Range("N8").Value = 16.485172
$ws.Range("O8").Value = 0.8161989011161211
$ws.Range("P8").Value = 0.8403205285996808
$ws.Range("Q8").Value = 53.85081453886932
$ws.Range("R8").Value = 323.1048872332159
$ws.Range("S8").Value = 0.5863223786081564
$ws.Range("T8").Value = 0.5291365773305833

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 9.799864
$ws.Range("H9").Value = 19.599728
$ws.Range("I9").Value = 0.7183572261692373
$ws.Range("J9").Value = 0.6296842208678898
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.657666
$ws.Range("N9").Value = 1.972998
$ws.Range("O9").Value = 0.09768528951377062
$ws.Range("P9").Value = 0.1005722428790014
$ws.Range("Q9").Value = 6.445037357423999
$ws.Range("R9").Value = 38.670224144544
$ws.Range("S9").Value = 0.07017293361265114
$ws.Range("T9").Value = 0.06332875439820017

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.799864
$ws.Range("H10").Value = 19.599728
$ws.Range("I10").Value = 0.7183572261692373
$ws.Range("J10").Value = 0.6296842208678898
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5797745
$ws.Range("N10").Value = 1.159549
$ws.Range("O10").Value = 0.08611580937010824
$ws.Range("P10").Value = 0.0591072285213179
$ws.Range("Q10").Value = 5.681711250667999
$ws.Range("R10").Value = 22.726845002672
$ws.Range("S10").Value = 0.06186191394842976
$ws.Range("T10").Value = 0.03721888913910638

Write-Host "Done updating rows 2-10"